$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.012386560440063
$ws.Range("B1").Value = 2.0063157081604
$ws.Range("C1").Value = 3.569568395614624
$ws.Range("D1").Value = 2.495243549346924
$ws.Range("E1").Value = 0.3393259644508362
